$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "survey" sheet

# Remove the stray, content-less column F cells (F4/F5) left over from an
# earlier edit - the updated form no longer has anything in column F.
$ws.Range("F4").Clear()
$ws.Range("F5").Clear()

# Insert a blank row before the old row 9 ("begin repeat / srepeat / New
# Sample") so the barcode/labid rows can move above the repeat and a new
# "Status Timestamp" question can be inserted between them.
$ws.Rows.Item(9).Insert()

# New row 9: Sample Tracking ID barcode question (was part of the old
# repeat, now asked once before the repeat begins).
$ws.Range("A9").Value = "barcode"
$ws.Range("B9").Value = "stid"
$ws.Range("C9").Value = "Sample Tracking ID"
$ws.Range("D9").Value = "yes"
$ws.Range("E9").ClearContents()

# Row 10 (was "begin repeat"): becomes the Lab ID barcode question.
$ws.Range("A10").Value = "barcode"
$ws.Range("B10").Value = "labid"
$ws.Range("C10").Value = "Lab ID"
$ws.Range("D10").Value = "yes"
$ws.Range("E10").ClearContents()

# Row 11 (was the "stid" barcode question): becomes the new Status
# Timestamp dateTime question.
$ws.Range("A11").Value = "dateTime"
$ws.Range("B11").Value = "labtime"
$ws.Range("C11").Value = "Status Timestamp"
$ws.Range("D11").Value = "yes"
$ws.Range("E11").ClearContents()

# Row 12 (was the "labid" barcode question): becomes "begin repeat" with
# the updated label "New Sample Test", and no longer carries
# required/appearance values.
$ws.Range("A12").Value = "begin repeat"
$ws.Range("B12").Value = "srepeat"
$ws.Range("C12").Value = "New Sample Test"
$ws.Range("D12").ClearContents()
$ws.Range("E12").ClearContents()

# Rows 13-16 (select_one labstatus/labtest/labreject, end repeat) keep
# their previous content unchanged - the row Insert() already shifted them
# down from 12-15 to 13-16.

# The geopoint question (now row 17) is no longer required.
$ws.Range("D17").ClearContents()
